$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for new column G
$ws.Cells.Item(1,7).Value = "Impressumcheck"

# Populate "+" marks in column G for the rows that satisfy the new
# "Impressumcheck" criterion (mirrors the "+" used throughout columns B-F).
$plusRows = @(2, 3, 6, 9, 10, 11, 12, 13, 14, 15, 16, 17, 18, 19, 20, 21, 22, 26, 27, 30, 33, 35, 36, 38, 39, 40, 42, 44, 46, 47, 50, 53, 55, 58, 60, 62)
foreach ($r in $plusRows) {
    $cell = $ws.Cells.Item($r, 7)
    $cell.Value = "+"
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4108
}

# Row(s) that only receive the centered formatting with no value in column G
$emptyStyledRows = @(49)
foreach ($r in $emptyStyledRows) {
    $cell = $ws.Cells.Item($r, 7)
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4108
}

# Restore selection to match the authored workbook
$ws.Range("G5").Select() | Out-Null
